$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting rows 108:132 down to 109:133.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new record.
$ws.Range("A108").Value = 1
$ws.Range("B108").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C108").Value = "Arica y Parinacota"
$ws.Range("D108").Value = 44841
$ws.Range("E108").Value = 15
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100106
$ws.Range("H108").Value = "Oleaginosos"
$ws.Range("I108").Value = 100106002
$ws.Range("J108").Value = "Palta"
$ws.Range("K108").Value = "Hass"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 480
$ws.Range("N108").Value = 23000
$ws.Range("O108").Value = 24000
$ws.Range("P108").Value = 23500
$ws.Range("Q108").Value = "`$/bandeja 10 kilos"
$ws.Range("R108").Value = "Perú"
$ws.Range("S108").Value = 2350
$ws.Range("T108").Value = 10
